{"js": "// Replace the division-problem text in each table cell with the new\n// expression, per the commit's diff. Each old expression is unique in the\n// document, so we can safely locate it with Body.search() (exact, case-\n// sensitive substring match) and swap just the `w:t` text via\n// insertText(..., \"Replace\"), leaving run/paragraph formatting untouched.\nconst replacements = [\n  [\"95\u00f76=\", \"70\u00f72=\"],\n  [\"70\u00f74=\", \"40\u00f74=\"],\n  [\"20\u00f77=\", \"20\u00f74=\"],\n  [\"54\u00f74=\", \"29\u00f77=\"],\n  [\"21\u00f78=\", \"31\u00f74=\"],\n  [\"70\u00f73=\", \"33\u00f72=\"],\n  [\"10\u00f78=\", \"42\u00f74=\"],\n  [\"76\u00f73=\", \"81\u00f79=\"],\n  [\"27\u00f73=\", \"52\u00f74=\"],\n  [\"50\u00f75=\", \"63\u00f72=\"],\n  [\"32\u00f76=\", \"32\u00f77=\"],\n  [\"90\u00f78=\", \"72\u00f72=\"],\n  [\"66\u00f77=\", \"42\u00f78=\"],\n  [\"49\u00f74=\", \"70\u00f77=\"],\n  [\"57\u00f76=\", \"54\u00f78=\"],\n  [\"86\u00f76=\", \"37\u00f79=\"],\n  [\"31\u00f73=\", \"89\u00f76=\"],\n  [\"78\u00f75=\", \"43\u00f76=\"],\n  [\"93\u00f79=\", \"89\u00f77=\"],\n  [\"68\u00f72=\", \"18\u00f77=\"],\n  [\"65\u00f72=\", \"30\u00f77=\"],\n  [\"29\u00f79=\", \"81\u00f72=\"],\n  [\"51\u00f75=\", \"20\u00f79=\"],\n  [\"28\u00f75=\", \"91\u00f73=\"],\n  [\"62\u00f74=\", \"80\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the division-problem text in each table cell with the new\n# expression, per the commit's diff. Each old expression occurs exactly\n# once in the document, so Find/Replace against the whole document\n# content is unambiguous and only rewrites the `w:t` text, leaving\n# run/paragraph formatting untouched.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"95\u00f76=\", \"70\u00f72=\"),\n    @(\"70\u00f74=\", \"40\u00f74=\"),\n    @(\"20\u00f77=\", \"20\u00f74=\"),\n    @(\"54\u00f74=\", \"29\u00f77=\"),\n    @(\"21\u00f78=\", \"31\u00f74=\"),\n    @(\"70\u00f73=\", \"33\u00f72=\"),\n    @(\"10\u00f78=\", \"42\u00f74=\"),\n    @(\"76\u00f73=\", \"81\u00f79=\"),\n    @(\"27\u00f73=\", \"52\u00f74=\"),\n    @(\"50\u00f75=\", \"63\u00f72=\"),\n    @(\"32\u00f76=\", \"32\u00f77=\"),\n    @(\"90\u00f78=\", \"72\u00f72=\"),\n    @(\"66\u00f77=\", \"42\u00f78=\"),\n    @(\"49\u00f74=\", \"70\u00f77=\"),\n    @(\"57\u00f76=\", \"54\u00f78=\"),\n    @(\"86\u00f76=\", \"37\u00f79=\"),\n    @(\"31\u00f73=\", \"89\u00f76=\"),\n    @(\"78\u00f75=\", \"43\u00f76=\"),\n    @(\"93\u00f79=\", \"89\u00f77=\"),\n    @(\"68\u00f72=\", \"18\u00f77=\"),\n    @(\"65\u00f72=\", \"30\u00f77=\"),\n    @(\"29\u00f79=\", \"81\u00f72=\"),\n    @(\"51\u00f75=\", \"20\u00f79=\"),\n    @(\"28\u00f75=\", \"91\u00f73=\"),\n    @(\"62\u00f74=\", \"80\u00f72=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
